# Bump the cached "datetimeFigureOut" date placeholder text from
# 2022/12/11 to 2023/1/31 across the slide master and every slide
# layout (mirrors PowerPoint re-caching the auto date field on save).

$p = $ppt.ActivePresentation
$newDate = "2023/1/31"

function Update-DatePlaceholder($shapes, $label) {
    $cnt = $shapes.Count
    for ($i = 1; $i -le $cnt; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            $phType = $sh.PlaceholderFormat.Type
            if ($phType -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder) {
            if ($sh.HasTextFrame) {
                $sh.TextFrame.TextRange.Text = $newDate
                Write-Host ($label + " shape " + $i + " date updated")
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
$masterShapes = $master.Shapes
Update-DatePlaceholder $masterShapes "Master"

# Every slide layout under the master
$layouts = $master.CustomLayouts
$layoutCount = $layouts.Count
for ($li = 1; $li -le $layoutCount; $li++) {
    $lay = $layouts.Item($li)
    $layShapes = $lay.Shapes
    $label = "Layout " + $li
    Update-DatePlaceholder $layShapes $label
}

Write-Host "Date placeholders updated"
